$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Início")
Write-Host $ws.Name
